# Manual xlsx fix for quick pipeline creation
#
# - venue_latitude (AK2) / venue_longitude (AL2): restore full floating
#   point precision for the coordinates instead of the truncated values.
# - predict_proba (AO2): the cell held a stray Python exception string
#   ("Wrong number of items passed 2, placement implies 1") that leaked
#   into the sheet instead of the actual predicted probability; replace
#   it with the numeric prediction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AK2").Value = 42.354148100000003
$ws.Range("AL2").Value = -71.104948899999997
$ws.Range("AO2").Value = 0.4

# Leave the selection where the author last left it.
$ws.Range("AO3").Select()
